$d = $word.ActiveDocument
$d.Content.Find.Execute("2609 Evergreen rd. Odenton, Maryland", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Odenton, Maryland", 2)
